$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values, forcing text format to preserve exact
# string representation (e.g. trailing zeros, dotted thousand separators).
$dCells = @("D2","D3","D5","D6","D7","D9","D10","D11","D13","D14","D15","D16","D18","D19","D20","D21","D22","D23","D25","D26","D27","D28","D29","D31","D32","D35","D37","D42","D44","D45","D46","D48","D49","D50")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "43.066.51"
$ws.Range("D3").Value = "2.305.86"
$ws.Range("D5").Value = "300.64"
$ws.Range("D6").Value = "97.95"
$ws.Range("D7").Value = "0.520"
$ws.Range("D9").Value = "0.520"
$ws.Range("D10").Value = "35.68"
$ws.Range("D11").Value = "0.0791"
$ws.Range("D13").Value = "17.90"
$ws.Range("D14").Value = "6.89"
$ws.Range("D15").Value = "2.664.18"
$ws.Range("D16").Value = "2.323.09"
$ws.Range("D18").Value = "42.978.13"
$ws.Range("D19").Value = "13.40"
$ws.Range("D20").Value = "0.0₃0908"
$ws.Range("D21").Value = "6.14"
$ws.Range("D22").Value = "68.36"
$ws.Range("D23").Value = "239.69"
$ws.Range("D25").Value = "1.00"
$ws.Range("D26").Value = "2.43"
$ws.Range("D27").Value = "24.72"
$ws.Range("D28").Value = "167.53"
$ws.Range("D29").Value = "9.18"
$ws.Range("D31").Value = "33.32"
$ws.Range("D32").Value = "5.22"
$ws.Range("D35").Value = "18.10"
$ws.Range("D37").Value = "0.0690"
$ws.Range("D42").Value = "2.008.60"
$ws.Range("D44").Value = "2.14"
$ws.Range("D45").Value = "10.08"
$ws.Range("D46").Value = "17.41"
$ws.Range("D48").Value = "54.40"
$ws.Range("D49").Value = "2.529.68"
$ws.Range("D50").Value = "73.92"

# Restore default (Normal) style so no stray number-format style is left
# on these cells, matching the original formatting.
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}

# Update Volume(1h) (column E) values.
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("E6").Value = "  -1.71%  "
$ws.Range("E7").Value = "  +3.91%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +1.89%  "
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("E19").Value = "  +7.59%  "
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("E23").Value = "  +1.74%  "
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  -0.70%  "
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("E30").Value = "  -11.10%  "
$ws.Range("E31").Value = "  -3.29%  "
$ws.Range("E32").Value = "  +4.81%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("E34").Value = "  +4.78%  "
$ws.Range("E35").Value = "  +4.71%  "
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("E39").Value = "  +1.05%  "
$ws.Range("E40").Value = "  +2.36%  "
$ws.Range("E41").Value = "  -2.37%  "
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("E44").Value = "  -2.95%  "
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("E47").Value = "  -1.82%  "
$ws.Range("E48").Value = "  -2.17%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("E50").Value = "  +6.00%  "
$ws.Range("E51").Value = "  +1.27%  "
